$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 520, shifting the existing rows 520:626 down to 521:627
$ws.Rows.Item(520).Insert()

# Populate the newly inserted row 520 with its data
$ws.Cells.Item(520, 1).Value2 = 4
$ws.Cells.Item(520, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(520, 3).Value2 = "Los Lagos"
$ws.Cells.Item(520, 4).Value2 = 45244
$ws.Cells.Item(520, 5).Value2 = 10
$ws.Cells.Item(520, 6).Value2 = 100112008
$ws.Cells.Item(520, 7).Value2 = "Coliflor"
$ws.Cells.Item(520, 8).Value2 = "Sin especificar"
$ws.Cells.Item(520, 9).Value2 = "Primera"
$ws.Cells.Item(520, 10).Value2 = 1500
$ws.Cells.Item(520, 11).Value2 = 1500
$ws.Cells.Item(520, 12).Value2 = 1500
$ws.Cells.Item(520, 13).Value2 = 1500
$ws.Cells.Item(520, 14).Value2 = "`$/unidad"
$ws.Cells.Item(520, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(520, 16).Value2 = 1500
$ws.Cells.Item(520, 17).Value2 = 1
$ws.Cells.Item(520, 18).Value2 = "Hortaliza"
